# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the "Latest Handoff
# Date"/"Latest Handoff Datetime" timestamps for the rows that were still
# pending (rows 7 and 10-16, i.e. everything at/after "f2b23054...") get
# refreshed to the new handoff run's timestamps, for all three sheets.

$wb = $excel.ActiveWorkbook

$overviewDate = "2016-35-18 05:35:15"
$zhcnDate     = "2016-03-18 05:35:12"
$dedeDate     = "2016-03-18 05:35:15"

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet: column D = "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = $overviewDate
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = $zhcnDate
}

# de-de sheet: column E = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = $dedeDate
}
